# Expense_Tracker.xlsx edit script
# - Re-orders several same-day Category/Amount pairs (rows 14-16, 52-55,
#   140-144, 149-151, 162-166, 188-190) to match the new canonical row order.
# - Inserts two brand-new test rows (2025/04/29 Restaurant $20, notes
#   "test"/"test3") ahead of the trailing 2025/04/30 rows, pushing those two
#   rows from 202/203 down to 204/205, and updates row 201's note to "test4".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Re-ordered Category/Amount pairs (date stays put; same-day rows were
#     shuffled) ---------------------------------------------------------
$ws.Cells.Item(14,2).Value = "Phone"
$ws.Cells.Item(14,3).Value = 105.95
$ws.Cells.Item(15,2).Value = "Barber"
$ws.Cells.Item(15,3).Value = 89.78
$ws.Cells.Item(16,2).Value = "Transportation"
$ws.Cells.Item(16,3).Value = 119.26

$ws.Cells.Item(52,3).Value = 169.37
$ws.Cells.Item(53,3).Value = 177.51
$ws.Cells.Item(54,2).Value = "Entertainment"
$ws.Cells.Item(54,3).Value = 91.19
$ws.Cells.Item(55,2).Value = "Restaurant"
$ws.Cells.Item(55,3).Value = 42.06

$ws.Cells.Item(140,2).Value = "Groceries"
$ws.Cells.Item(140,3).Value = 21.07
$ws.Cells.Item(142,2).Value = "Toters"
$ws.Cells.Item(142,3).Value = 169.93
$ws.Cells.Item(144,2).Value = "Restaurant"
$ws.Cells.Item(144,3).Value = 45.66

$ws.Cells.Item(149,2).Value = "Snacks"
$ws.Cells.Item(149,3).Value = 6.21
$ws.Cells.Item(150,2).Value = "Barber"
$ws.Cells.Item(150,3).Value = 199
$ws.Cells.Item(151,2).Value = "Entertainment"
$ws.Cells.Item(151,3).Value = 101.09

$ws.Cells.Item(162,2).Value = "Barber"
$ws.Cells.Item(162,3).Value = 180.03
$ws.Cells.Item(163,2).Value = "Restaurant"
$ws.Cells.Item(163,3).Value = 105.4
$ws.Cells.Item(164,2).Value = "Snacks"
$ws.Cells.Item(164,3).Value = 177.11
$ws.Cells.Item(165,2).Value = "Phone"
$ws.Cells.Item(165,3).Value = 159.19
$ws.Cells.Item(166,2).Value = "Groceries"
$ws.Cells.Item(166,3).Value = 115.05

$ws.Cells.Item(188,2).Value = "Entertainment"
$ws.Cells.Item(188,3).Value = 199.43
$ws.Cells.Item(189,2).Value = "Groceries"
$ws.Cells.Item(189,3).Value = 171.79
$ws.Cells.Item(190,2).Value = "Restaurant"
$ws.Cells.Item(190,3).Value = 183.49

# --- Tail of the sheet: update note on row 201, insert two brand-new rows
#     (202, 203), pushing the old 202/203 down to 204/205 ----------------
$ws.Cells.Item(201,4).Value = "test4"

$ws.Rows("202:203").Insert()

# Row 201 (2025/04/29, Restaurant, 20) already has the text-formatted date
# we need, so copy it down into the two freshly inserted rows instead of
# re-typing the date (typing "2025/04/29" directly would be auto-coerced
# to a real date serial instead of staying plain text).
$ws.Rows(201).Copy()
$ws.Rows(202).PasteSpecial()
$ws.Rows(201).Copy()
$ws.Rows(203).PasteSpecial()

$ws.Cells.Item(202,4).Value = "test"
$ws.Cells.Item(203,4).Value = "test3"

Write-Output ("Final UsedRange=" + $ws.UsedRange.Address())
